$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# New shared string for the two new rows
$newTask = "Unsuccessful, reverted attempt to optimize context switch code"

# Row 19: 24 June 2013 (serial 41449), 2h, task description
$ws.Range("A18").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 41449
$ws.Range("B19").Value = 2
$ws.Range("D19").Value = $newTask

# Row 20: 25 June 2013 (serial 41450), 1h, same task description
$ws.Range("A18").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 41450
$ws.Range("B20").Value = 1
$ws.Range("D20").Value = $newTask

# Update selection to match the committed state
$ws.Range("A19").Select()
